$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.223.64'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.894.59'
$ws.Range('E3').Value = '  +2.26%  '
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.41'
$ws.Range('E5').Value = '  +2.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.657'
$ws.Range('E6').Value = '  +6.03%  '
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.37'
$ws.Range('E8').Value = '  -1.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.343'
$ws.Range('E9').Value = '  +5.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '50.55'
$ws.Range('E10').Value = '  +8.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0708'
$ws.Range('E11').Value = '  +2.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0997'
$ws.Range('E12').Value = '  +0.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.170.63'
$ws.Range('E13').Value = '  +2.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.95'
$ws.Range('E14').Value = '  +5.31%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.912.81'
$ws.Range('E15').Value = '  +3.15%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.691'
$ws.Range('E16').Value = '  +2.81%  '
$ws.Range('E17').Value = '  +1.81%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.201.10'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.05'
$ws.Range('E19').Value = '  +1.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0812'
$ws.Range('E20').Value = '  +2.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '240.43'
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.39'
$ws.Range('E22').Value = '  +2.15%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.39'
$ws.Range('E25').Value = '  +31.07%  '
$ws.Range('E26').Value = '  +0.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '169.76'
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('E28').Value = '  +4.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.14'
$ws.Range('E29').Value = '  +3.44%  '
$ws.Range('E30').Value = '  +2.51%  '
$ws.Range('E31').Value = '  +3.29%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0559'
$ws.Range('E32').Value = '  +1.34%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.937'
$ws.Range('E33').Value = '  +18.92%  '
$ws.Range('E34').Value = '  -0.35%  '
$ws.Range('E35').Value = '  +2.91%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.02'
$ws.Range('E37').Value = '  +1.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.32'
$ws.Range('E38').Value = '  +2.68%  '
$ws.Range('E39').Value = '  +1.80%  '
$ws.Range('E40').Value = '  +3.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0636'
$ws.Range('E41').Value = '  +14.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.96'
$ws.Range('E42').Value = '  +8.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '88.82'
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.336.25'
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '48.95'
$ws.Range('E45').Value = '  +43.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.35'
$ws.Range('E46').Value = '  +2.76%  '
$ws.Range('E47').Value = '  -1.51%  '
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.50'
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.078.32'
$ws.Range('E50').Value = '  +2.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.20'
$ws.Range('E51').Value = '  -14.56%  '
